$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de), row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Locale detail sheets: "Status" column (C), row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Shrink "Status" related columns' width ---
# Overview sheet columns E:F (zh-cn / de-de status columns)
$wsOverview.Columns("E:F").ColumnWidth = 12.576851254417766

# zh-cn / de-de sheets column C ("Status")
$wsZhCn.Columns("C:C").ColumnWidth = 12.576851254417766
$wsDeDe.Columns("C:C").ColumnWidth = 12.576851254417766
